$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Medellin -> Cali, 2 adults, 0 kids, 1 infant
$ws.Range("A3").Value = "Medellin"
$ws.Range("B3").Value = "Cali"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

# Row 4: Bogota -> Cucuta, 1 adult, 1 kid, 0 infants
$ws.Range("A4").Value = "Bogota"
$ws.Range("B4").Value = "Cucuta"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

# Update selection to match target state
$ws.Range("F4").Select()

$wb.Save()
